function Trim-ParaText($s) {
    # Paragraph.Range.Text carries a trailing paragraph/cell mark
    # (chr 13, and chr 7 inside table cells) - strip it for comparisons.
    return $s.TrimEnd([char]13, [char]7)
}

$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) Contact line: fix the "Conntact" typo and double the {{ }}
#    braces around the Contact.Cellular / Contact.Email placeholders.
# ------------------------------------------------------------------
$d.Content.Find.Execute("Conntact.Cellular", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "Contact.Cellular", 2) | Out-Null

$d.Content.Find.Execute("{Contact.Cellular}  | {Contact.Email}", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "{{Contact.Cellular}}  | {{Contact.Email}}", 2) | Out-Null

# ------------------------------------------------------------------
# 2) Remove the "Dear," heading paragraph entirely (its paragraph
#    mark goes with it, leaving the blank paragraph after it intact).
# ------------------------------------------------------------------
foreach ($p in $d.Paragraphs) {
    if ((Trim-ParaText $p.Range.Text) -eq "Dear,") {
        $p.Range.Delete()
        break
    }
}

# ------------------------------------------------------------------
# 3) Double the {{ }} braces around CoverLetter. It is alone in its
#    own paragraph, so inserting right at the paragraph boundaries
#    gives the new brace its own run either side of the placeholder.
# ------------------------------------------------------------------
foreach ($p in $d.Paragraphs) {
    if ((Trim-ParaText $p.Range.Text) -eq "{CoverLetter}") {
        $p.Range.InsertBefore("{")
        break
    }
}
foreach ($p in $d.Paragraphs) {
    if ((Trim-ParaText $p.Range.Text) -eq "{{CoverLetter}") {
        $p.Range.InsertAfter("}")
        break
    }
}

# ------------------------------------------------------------------
# 4) Collapse the "Sincerely," / "{{Name}}" sign-off paragraphs into
#    a single paragraph containing "=".
# ------------------------------------------------------------------
foreach ($p in $d.Paragraphs) {
    if ((Trim-ParaText $p.Range.Text) -eq "Sincerely,") {
        $p.Range.Text = "="
        break
    }
}
foreach ($p in $d.Paragraphs) {
    if ((Trim-ParaText $p.Range.Text) -eq "{{Name}}" -and $p.Range.Start -gt 0) {
        $p.Range.Delete()
        break
    }
}
